$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4090.9092

$ws.Range("H86").Value = 150068430
$ws.Range("I86").Value = 90911300
$ws.Range("J86").Value = 312750560
$ws.Range("K86").Value = 90911300
$ws.Range("L86").Value = 312750560
$ws.Range("M86").Value = -90910177
$ws.Range("N86").Value = -312752806

$ws.Range("H89").Value = 150068430
$ws.Range("I89").Value = 90911300
$ws.Range("J89").Value = 312750560
$ws.Range("K89").Value = 454556500
$ws.Range("L89").Value = 1563752800
$ws.Range("M89").Value = -454550884
$ws.Range("N89").Value = -1563764032

$ws.Range("H92").Value = 41667710
$ws.Range("I92").Value = 45455680
$ws.Range("K92").Value = 45455680
$ws.Range("M92").Value = -45454432

$ws.Range("H98").Value = 1501.8695
$ws.Range("I98").Value = 1387.762
$ws.Range("K98").Value = 1387.762
$ws.Range("M98").Value = 110.2380000000001

$ws.Range("H107").Value = 12490.4
$ws.Range("I107").Value = 12490.4
$ws.Range("K107").Value = 12490.4
$ws.Range("M107").Value = -10570.4

$ws.Range("H122").Value = 1501.8695
$ws.Range("I122").Value = 1387.762
$ws.Range("K122").Value = 4163.286
$ws.Range("M122").Value = -1713.286

$ws.Range("H132").Value = 3749.8113
$ws.Range("J132").Value = 2939.5
$ws.Range("L132").Value = 8818.5
$ws.Range("N132").Value = -13878.5

$ws.Range("H135").Value = 1259.762
$ws.Range("I135").Value = 1151.85
$ws.Range("K135").Value = 10366.65
$ws.Range("M135").Value = -7831.65

$ws.Range("H138").Value = 1958.6938
$ws.Range("I138").Value = 1186.375
$ws.Range("J138").Value = 2491.3276
$ws.Range("K138").Value = 3559.125
$ws.Range("L138").Value = 7473.9828
$ws.Range("M138").Value = 1580.875
$ws.Range("N138").Value = -17753.9828

$ws.Range("H141").Value = 1318.8235
$ws.Range("I141").Value = 1214.8667
$ws.Range("K141").Value = 3644.6001
$ws.Range("M141").Value = 1535.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17757440
$ws.Range("I32").Value = 18384040
$ws.Range("K32").Value = 18384040
$ws.Range("M32").Value = -18383753

$ws.Range("H45").Value = 3732.4
$ws.Range("I45").Value = 3785.0952
$ws.Range("J45").Value = 3609.4443
$ws.Range("K45").Value = 3785.0952
$ws.Range("L45").Value = 3609.4443
$ws.Range("M45").Value = -3408.0952
$ws.Range("N45").Value = -4363.4443

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1350060.6
$ws.Range("I134").Value = 1554900.5
$ws.Range("J134").Value = 3970
$ws.Range("K134").Value = 4664701.5
$ws.Range("L134").Value = 11910
$ws.Range("M134").Value = -4662166.5
$ws.Range("N134").Value = -16980

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 3032
$ws.Range("I5").Value = 65
$ws.Range("J5").Value = 5999
$ws.Range("K5").Value = 65
$ws.Range("L5").Value = 5999
$ws.Range("M5").Value = 47
$ws.Range("N5").Value = -6223

$ws.Range("H15").Value = 1006
$ws.Range("I15").Value = 1006
$ws.Range("K15").Value = 1006
$ws.Range("M15").Value = -836

$ws.Range("H31").Value = 3778.577
$ws.Range("I31").Value = 1616.3611
$ws.Range("J31").Value = 8643.5625
$ws.Range("K31").Value = 1616.3611
$ws.Range("L31").Value = 8643.5625
$ws.Range("M31").Value = -1321.3611
$ws.Range("N31").Value = -9233.5625

$ws.Range("H34").Value = 3778.577
$ws.Range("I34").Value = 1616.3611
$ws.Range("J34").Value = 8643.5625
$ws.Range("K34").Value = 1616.3611
$ws.Range("L34").Value = 8643.5625
$ws.Range("M34").Value = -1414.3611
$ws.Range("N34").Value = -9047.5625

$ws.Range("H122").Value = 3574286.8
$ws.Range("I122").Value = 4169545
$ws.Range("K122").Value = 12508635
$ws.Range("M122").Value = -12506185

$ws.Range("H132").Value = 2955.883
$ws.Range("I132").Value = 2982.7627
$ws.Range("K132").Value = 8948.288100000002
$ws.Range("M132").Value = -6418.288100000002

$ws.Range("H134").Value = 2292.0278
$ws.Range("I134").Value = 2087.742
$ws.Range("K134").Value = 6263.226000000001
$ws.Range("M134").Value = -3728.226000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9997.666999999999
$ws.Range("I3").Value = 9997
$ws.Range("K3").Value = 29991
$ws.Range("M3").Value = -29879

$ws.Range("H8").Value = 400
$ws.Range("I8").Value = 400
$ws.Range("K8").Value = 1200
$ws.Range("M8").Value = -1061

$ws.Range("H14").Value = 596.1818
$ws.Range("I14").Value = 596.1818
$ws.Range("K14").Value = 1788.5454
$ws.Range("M14").Value = -1615.5454

$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9470

$ws.Range("H75").Value = 497
$ws.Range("J75").Value = 496.33334
$ws.Range("L75").Value = 1489.00002
$ws.Range("N75").Value = -3485.00002

$ws.Range("H78").Value = 497
$ws.Range("J78").Value = 496.33334
$ws.Range("L78").Value = 4467.00006
$ws.Range("N78").Value = -14451.00006

$ws.Range("H97").Value = 351
$ws.Range("J97").Value = 230.16667
$ws.Range("L97").Value = 690.50001
$ws.Range("N97").Value = -1682.50001

$ws.Range("H131").Value = 1585.9474
$ws.Range("I131").Value = 677
$ws.Range("J131").Value = 2116.1667
$ws.Range("K131").Value = 2031
$ws.Range("L131").Value = 6348.500100000001
$ws.Range("M131").Value = 3009
$ws.Range("N131").Value = -16428.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2200
$ws.Range("I80").Value = 2200
$ws.Range("J80").Value = 2200
$ws.Range("K80").Value = 2200
$ws.Range("L80").Value = 2200
$ws.Range("M80").Value = -1202
$ws.Range("N80").Value = -4196

$ws.Range("H83").Value = 2200
$ws.Range("I83").Value = 2200
$ws.Range("J83").Value = 2200
$ws.Range("K83").Value = 11000
$ws.Range("L83").Value = 11000
$ws.Range("M83").Value = -6008
$ws.Range("N83").Value = -20984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 68603.42
$ws.Range("I132").Value = 90987.22
$ws.Range("K132").Value = 272961.66
$ws.Range("M132").Value = -270431.66

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 9188.691999999999
$ws.Range("I96").Value = 2914.4285
$ws.Range("K96").Value = 2914.4285
$ws.Range("M96").Value = -1541.4285

$ws.Range("H126").Value = 4668.6665
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4668.6665
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 14005.9995
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -18945.9995

$ws.Range("H132").Value = 4125.8423
$ws.Range("I132").Value = 4226.4546
$ws.Range("J132").Value = 3987.5
$ws.Range("K132").Value = 12679.3638
$ws.Range("L132").Value = 11962.5
$ws.Range("M132").Value = -10149.3638
$ws.Range("N132").Value = -17022.5

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 42968.68
$ws.Range("I136").Value = 2637.111
$ws.Range("J136").Value = 146678.42
$ws.Range("K136").Value = 7911.333
$ws.Range("L136").Value = 440035.26
$ws.Range("M136").Value = -5361.333
$ws.Range("N136").Value = -445135.26
